$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6711
    $ws.Range("F3").Value = 45
    $ws.Range("F5").Value = 1051
    $ws.Range("F6").Value = 142
}
